# Update the stale Python object memory addresses embedded in the
# repr() strings stored in A3, C3, A4 and C4 (NamedFeatureSelector
# instances re-created on a later run get a new id()/address).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a3 = $ws.Range("A3").Value()
$ws.Range("A3").Value = $a3.Replace("0x7fcd83514160", "0x7f912432ae50")

$c3 = $ws.Range("C3").Value()
$ws.Range("C3").Value = $c3.Replace("0x7fccd04b2190", "0x7f91044a5460")

$a4 = $ws.Range("A4").Value()
$ws.Range("A4").Value = $a4.Replace("0x7fcd830ec7c0", "0x7f9104649220")

$c4 = $ws.Range("C4").Value()
$ws.Range("C4").Value = $c4.Replace("0x7fcd83162fa0", "0x7f9104555eb0")
